$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Email" header in C1, matching the header format used in A1/B1 ---
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C1").Value = "Email"

# --- Add the e-mail address in C2 and turn it into a mailto hyperlink ---
$ws.Range("C2").Value = "michael.aruebo@cyberspace.net.ng"
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:michael.aruebo@cyberspace.net.ng")

# --- Resize columns B and C for the new layout ---
$ws.Columns("B").ColumnWidth = 26.5
$ws.Columns("C").ColumnWidth = 31.666666666666668

# --- Update the active selection saved in the sheet view ---
$ws.Range("C6").Select() | Out-Null

Write-Output "applied staff-name sample email column update"
